# "Check max time for employee"
# Expands the single-employee-category WorkTime calculation into a
# multi-category comparison (0x0 .. 5x8 shift patterns) for "Средний
# медперсонал" at a 19-day norm, and raises the reference day norm
# from 72,36 to 148,12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: parameters used by every WorkTime() call ------------------
$ws.Range("A1").Value = 19
$ws.Range("B1").Value = "1,0=148,12"
$ws.Range("C1").Value = "Средний медперсонал"

# ---- Row 3: the six shift patterns being compared ----------------------
$ws.Range("A3").Value = "0x0"
$ws.Range("B3").Value = "1x9"
$ws.Range("C3").Value = "2x9"
$ws.Range("D3").Value = "3x8"
$ws.Range("E3").Value = "4x5"
$ws.Range("F3").Value = "5x8"

# ---- Row 5: second line of each two-row spilled CSE array result -------
# Must be written BEFORE the array formula is entered on row 4, otherwise
# the cell becomes part of the array block and can no longer be edited
# on its own ("You cannot change part of an array").
$ws.Range("B5").Value = "1x7,30"
$ws.Range("C5").Value = "1x6,54"
$ws.Range("D5").Value = "1x6,36"
$ws.Range("E5").Value = "1x7,0"
$ws.Range("F5").Value = "1x5,54"

# ---- Row 4: legacy CSE array formulas calling the WorkTime add-in ------
$ws.Range("A4").FormulaArray = '=_xll.WorkTime($B1,A3,$C1,$A1)'
$ws.Range("B4:B5").FormulaArray = '=_xll.WorkTime($B1,B3,$C1,$A1)'
$ws.Range("C4:C5").FormulaArray = '=_xll.WorkTime($B1,C3,$C1,$A1)'
$ws.Range("D4:D5").FormulaArray = '=_xll.WorkTime($B1,D3,$C1,$A1)'
$ws.Range("E4:E5").FormulaArray = '=_xll.WorkTime($B1,E3,$C1,$A1)'
$ws.Range("F4:F5").FormulaArray = '=_xll.WorkTime($B1,F3,$C1,$A1)'

# ---- Cosmetics: column width, selection ---------------------------------
$ws.Range("C1").ColumnWidth = 9.59
$ws.Range("C4").Select()
